$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.659.73"
$ws.Range("E2").Value = "  -2.09%  "
$ws.Range("D3").Value = "2.890.21"
$ws.Range("E3").Value = "  -2.32%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "568.59"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.50%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.98"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.74%  "
$ws.Range("E7").Value = "  +0.17%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.506"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.68%  "
$ws.Range("D9").Value = "2.890.56"
$ws.Range("E9").Value = "  -2.23%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.63"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -8.59%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.147"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.69%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.434"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.43%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000231"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.97%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "31.83"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.24%  "
$ws.Range("E15").Value = "  -0.69%  "
$ws.Range("D16").Value = "3.373.79"
$ws.Range("E16").Value = "  -2.23%  "
$ws.Range("D17").Value = "61.669.91"
$ws.Range("E17").Value = "  -2.01%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.61"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.14%  "
$ws.Range("D19").Value = "2.885.11"
$ws.Range("E19").Value = "  -2.62%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "432.61"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.04%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.17"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.41%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.654"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.38%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.88"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.37%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "79.52"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.76%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.81"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.73%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.19"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -9.46%  "
$ws.Range("E27").Value = "  -0.09%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.02"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -6.09%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0000102"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.95%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.03%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.50"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.18%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.05"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -6.20%  "
$ws.Range("B33").Value = "FirstDigitalUSD"
$ws.Range("C33").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.00"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.03%  "
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.106"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.83%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "25.52"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.960"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.23%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.41"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.11%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "48.96"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.43%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.93"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -6.56%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.78"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -11.65%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.115"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.00%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.23"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.65%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "39.17"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.93%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.267"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.49%  "
$ws.Range("D45").Value = "2.695.49"
$ws.Range("E45").Value = "  -0.75%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "132.80"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.60%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0333"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.08%  "
$ws.Range("E48").Value = "  +0.04%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "337.44"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -7.37%  "
$ws.Range("E50").Value = "  -1.96%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "21.56"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -6.05%  "
